$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.732.19"
$ws.Range("E2").Value = "  +1.96%  "
# Row 3
$ws.Range("D3").Value = "1.924.23"
$ws.Range("E3").Value = "  +0.68%  "
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9991"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -1.13%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "335.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.88%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.99%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4672"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.30%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4147"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.29%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.52"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.14%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08066"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.80%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.022"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.86%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.33"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.22%  "
# Row 13
$ws.Range("D13").Value = "1.922.51"
$ws.Range("E13").Value = "  -1.22%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.024"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.17%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.212"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.02%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "90.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.31%  "
# Row 17
$ws.Range("E17").Value = "  -0.96%  "
# Row 18
$ws.Range("E18").Value = "  -0.68%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06593"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.90%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.58%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.65%  "
# Row 22
$ws.Range("D22").Value = "29.673.89"
$ws.Range("E22").Value = "  +1.71%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.568"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.06%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.63"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.58%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.201"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.88%  "
# Row 26
$ws.Range("D26").Value = "2.153.02"
$ws.Range("E26").Value = "  -0.76%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.66%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.96"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.52%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.158"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.69%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.770"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.62%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "118.03"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.40%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.049"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.59%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09473"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.87%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.438"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.06%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.447"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.19%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.527"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.44%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06163"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.15%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02272"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.44%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.480"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.22%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.180"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.46%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5921"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.59%  "
# Row 42
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1850"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.73%  "
# Row 43
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.29"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.47%  "
# Row 44
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.362"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.46%  "
# Row 45
$ws.Range("B45").Value = "WEMIXTOKEN"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.256"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.23%  "
# Row 46
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07530"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.99%  "
# Row 47
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5606"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.13%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "12.20"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.62%  "
# Row 49
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.944"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.14%  "
# Row 50
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "113.21"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.21%  "
# Row 51
$ws.Range("B51").Value = "WOONetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3019"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +9.57%  "
